$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# --- Row 7 ---
# A7: date, reuse style from A2 (numFmtId 14)
$ws.Range("A2").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value = 40220

# B7: from-time, reuse style from B2 (numFmtId 20)
$ws.Range("B2").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("B7").Value = 0.4375

# C7: till-time, reuse style from C2 (numFmtId 20)
$ws.Range("C2").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("C7").Value = 0.60416666666666663

# D7: plain number, default style
$ws.Range("D7").Value = 4

# E7 / F7: shared strings
$ws.Range("E7").Value = "Design"
$ws.Range("F7").Value = "MessageQueue and rendering design, Cutting out tiles"

# --- Row 8 ---
# A8: date, reuse style from A2 (numFmtId 14)
$ws.Range("A2").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = 40225

# B8: from-time, reuse style from B2 (numFmtId 20)
$ws.Range("B2").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("B8").Value = 0.45833333333333331

# C8: till-time, reuse style from C2 (numFmtId 20)
$ws.Range("C2").Copy()
$ws.Range("C8").PasteSpecial(-4122)
$ws.Range("C8").Value = 0.60416666666666663

# D8: text "3.5" stored as shared string with default style.
# Build it as a formula returning text, then collapse to a value in place
# so no new number-format style gets allocated.
$ws.Range("D8").Formula = '="3.5"'
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)

# E8 / F8: shared strings
$ws.Range("E8").Value = "Analysis"
$ws.Range("F8").Value = "Plan of Attack"

$excel.CutCopyMode = 0

# Update the saved selection to match the author's final cursor position
$ws.Range("F17").Select()
